# Apply the ENTSOE "Actual Production Wind" data refresh:
#  - Column A (timestamps) for rows 2-97 are shifted forward by 4 days
#    (45821.x -> 45825.x, 45822 -> 45826)
#  - Column B (Actual Production MW) for rows 2-41 gets new values from the
#    latest data pull; rows 42-97 remain 0 (unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Actual Production (MW) values for rows 2 through 41 (rows 42-97 stay 0)
$newB = @(186,163,158,157,172,180,179,169,150,134,132,139,154,169,180,216,217,238,240,247,255,276,298,307,327,363,342,302,269,250,255,260,262,264,277,293,331,345,367,373)

for ($r = 2; $r -le 97; $r++) {
    # Shift the timestamp in column A forward by 4 days (keeps the same time-of-day)
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value2 = $cellA.Value2 + 4

    # Update column B where a new value was provided by this data pull
    $idx = $r - 2
    if ($idx -lt $newB.Length) {
        $ws.Cells.Item($r, 2).Value2 = $newB[$idx]
    }
}
